$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "305.66"
Set-TextValue "E2" "1.53%"
Set-TextValue "D3" "31.88"
Set-TextValue "E3" "0.91%"
Set-TextValue "D4" "5.217"
Set-TextValue "E4" "1.03%"
Set-TextValue "D5" "0.07476"
Set-TextValue "E5" "-0.24%"
Set-TextValue "D6" "2.274"
Set-TextValue "E6" "34.60%"
Set-TextValue "D7" "7.980"
Set-TextValue "E7" "1.56%"
Set-TextValue "D8" "3.873"
Set-TextValue "E8" "1.02%"
Set-TextValue "D9" "0.9229"
Set-TextValue "E9" "0.01%"
Set-TextValue "D10" "0.1740"
Set-TextValue "E10" "1.55%"
Set-TextValue "D11" "0.07693"
Set-TextValue "E11" "0.80%"
Set-TextValue "D12" "0.08207"
Set-TextValue "E12" "2.20%"
Set-TextValue "D13" "0.02990"
Set-TextValue "E13" "0.05%"
Set-TextValue "D14" "0.09959"
Set-TextValue "E14" "0.76%"
Set-TextValue "D15" "0.001507"
Set-TextValue "E15" "0.89%"
Set-TextValue "D16" "0.006101"
Set-TextValue "E16" "-1.58%"
Set-TextValue "D17" "3.491"
Set-TextValue "E17" "1.04%"
Set-TextValue "E18" "0.04%"
Set-TextValue "D19" "0.3266"
Set-TextValue "E19" "-0.84%"
Set-TextValue "D20" "0.1334"
Set-TextValue "E20" "-1.09%"
Set-TextValue "D21" "4.637"
Set-TextValue "E21" "1.33%"
Set-TextValue "D22" "0.04603"
Set-TextValue "E22" "-1.40%"
Set-TextValue "D23" "0.1559"
Set-TextValue "E23" "0.50%"
Set-TextValue "D24" "0.001235"
Set-TextValue "E24" "0.99%"
Set-TextValue "D25" "0.004534"
Set-TextValue "E25" "2.33%"
Set-TextValue "D26" "0.0001294"
Set-TextValue "E26" "-7.62%"
Set-TextValue "D27" "0.0002730"
Set-TextValue "E27" "51.92%"
Set-TextValue "D39" "0.01799"
Set-TextValue "E39" "7.99%"
Set-TextValue "D40" "0.04571"
Set-TextValue "E40" "0.53%"
Set-TextValue "D41" "0.007362"
Set-TextValue "E41" "5.62%"
Set-TextValue "D42" "0.1368"
Set-TextValue "E42" "1.90%"
Set-TextValue "D43" "0.002110"
Set-TextValue "E43" "2.38%"
Set-TextValue "D44" "0.01080"
Set-TextValue "E44" "-12.57%"
Set-TextValue "D45" "0.00006435"
Set-TextValue "E45" "6.55%"
Set-TextValue "E46" "-57.20%"
Set-TextValue "D47" "0.009857"
Set-TextValue "E47" "-19.52%"

Write-Output "Done applying changes"